$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the timestamp stored in A95 (row for 2024-07-05) ---
$ws.Range("A95").Value = 45478.2916666667

# --- Append new row 96 (2024-07-08 trading data) ---

# Column A: datetime value, re-using A95's date/time number format (style s="1")
$ws.Range("A96").Value = 45481.6314351852
$ws.Range("A95").Copy() | Out-Null
$ws.Range("A96").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B96").Value = 4500
$ws.Range("C96").Value = 3.4300000667572
$ws.Range("D96").Value = 3.3199999332428
$ws.Range("E96").Value = 3.4300000667572
$ws.Range("F96").Value = 3.35999989509583

# Column G: text representation of the close value (stored as a shared string)
$ws.Range("G96").Formula = '="3.35999989509583"'
$ws.Range("G96").Copy() | Out-Null
$ws.Range("G96").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$ws.Range("H96").Value = "ESPE.MI"

$excel.CutCopyMode = 0
